$d = $word.ActiveDocument

# --- Item 1: Reconocimiento y reconocimiento (bold heading) / trailing text ---
$d.Content.Find.Execute(
    "Reconocimiento y reconocimiento", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Concienciación y reconocimiento de marca limitados",
    2)

$d.Content.Find.Execute(
    " de marca limitados: lograr visibilidad en estos nuevos mercados es un obstáculo principal, lo que requiere un sólido esfuerzo de marketing para crear la presencia de marca de Adatum desde cero.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    ": un obstáculo principal es lograr visibilidad en estos nuevos mercados, lo que requiere un sólido esfuerzo de marketing para crear la presencia de marca de Adatum desde cero.",
    2)

# --- Item 2: Intensa competencia (bold heading) / trailing text ---
$d.Content.Find.Execute(
    "Intensa competencia", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Competencia intensa",
    2)

$d.Content.Find.Execute(
    ": el sector de servicios en la nube en Canadá es ferozmente competitivo, con numerosos jugadores.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    ": el sector de servicios en la nube en Canadá es ferozmente competitivo, con numerosos participantes.",
    2)

# --- Item 3: Diversas preferencias y expectativas (bold heading) / trailing text ---
$d.Content.Find.Execute(
    "Diversas preferencias y expectativas", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Preferencias y expectativas diversas de los clientes",
    2)

$d.Content.Find.Execute(
    " de los clientes: adaptar productos y marketing para alinearse con las diversas demandas de estos mercados es fundamental para replicar con empresas y consumidores locales.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    ": adaptar los productos y el marketing para alinearse con las diversas demandas de estos mercados es fundamental para resonar con las empresas y los consumidores locales.",
    2)

# --- Item 4: Desafíos (bold heading) / trailing text ---
$d.Content.Find.Execute(
    "Desafíos", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Desafíos normativos y de cumplimiento",
    2)

$d.Content.Find.Execute(
    " normativos y de cumplimiento: Adatum se enfrenta a la compleja tarea de navegar por la privacidad, la seguridad y las regulaciones operativas de la región, lo que necesita esfuerzos de cumplimiento diligentes.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    ": Adatum se enfrenta a la compleja tarea de navegar por la privacidad, la seguridad y las regulaciones operativas de la región, requiriendo esfuerzos de cumplimiento diligentes.",
    2)

# --- Item 5: Complejidad (bold heading) / trailing text ---
$d.Content.Find.Execute(
    "Complejidad", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Complejidad operativa y logística",
    2)

$d.Content.Find.Execute(
    " operativa y logística: el establecimiento de operaciones eficientes entre regiones presenta desafíos logísticos, especialmente en el mantenimiento de altos niveles de servicio y la administración de centros de datos en ubicaciones geográficas.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    ": el establecimiento de operaciones eficientes entre regiones presenta desafíos logísticos, especialmente en el mantenimiento de altos niveles de servicio y la administración de centros de datos en todas las ubicaciones geográficas.",
    2)
